$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22 (A22 = "Menton")
$ws.Range("B22").Value = 12
$ws.Range("C22").Value = 6.911295917190489
$ws.Range("D22").Value = 4.145507740228843
$ws.Range("F22").Value = 3.824519598430504
$ws.Range("G22").Value = 6.103572369942697
$ws.Range("H22").Value = 7.364563152662695

# Row 28 (A28 = "Pogonion")
$ws.Range("B28").Value = 12
$ws.Range("C28").Value = 7.24866125064084
$ws.Range("D28").Value = 4.534872458372755
$ws.Range("F28").Value = 3.869574785684582
$ws.Range("G28").Value = 6.49753587688329
$ws.Range("H28").Value = 8.701017672821747
